# Update Rev 2 Files
# - Row 17: Designator list loses "R18"; Quantity drops from 7 to 6.
# - Row 32 (U9): part swapped from Bosch Sensortec BMI088 / BMI055 footprint
#   to TDK InvenSense ICM-42688-P, with LibRef/Mfr Part Desc/Supplier Part
#   Number updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a text value into a cell while preserving the cell's
# existing style/number-format (a plain .Value assignment would either
# coerce numeric-looking strings into numbers, or strip the quotePrefix
# style bit off alphabetic strings). We stage the text in a scratch cell
# (forced to Text format so digit strings remain text), copy it, and use
# PasteSpecial Values-only so only the cell content changes.
function Set-TextValue {
    param($targetCell, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# --- Row 17 -----------------------------------------------------------
Set-TextValue $ws.Range("B17") "R4, R6, R10, R13, R15, R17"
$ws.Range("C17").Value = 6

# --- Row 32 (U9) --------------------------------------------------------
Set-TextValue $ws.Range("A32") "31021"
Set-TextValue $ws.Range("D32") "TDK InvenSense"
Set-TextValue $ws.Range("E32") "ICM-42688-P"
Set-TextValue $ws.Range("F32") "ICM-42688-P"
Set-TextValue $ws.Range("G32") "MOTION SENSOR"
Set-TextValue $ws.Range("H32") "ICM-42688-P"
Set-TextValue $ws.Range("J32") "1428-ICM-42688-PTR-ND"
